$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting rows 10-46 down to 11-47
$ws.Rows("10:10").Insert()

# Fill in the new row 10 with its data
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44910
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 300000000
$ws.Range("G10").Value = "Espárragos"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 1200
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 950
$ws.Range("N10").Value = "$/kilo"
$ws.Range("O10").Value = "Provincia de Diguillín"
$ws.Range("P10").Value = 950
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
